$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$newRows = @(
  @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2021.xlsx", "2021"),
  @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2022.xlsx", "2022")
)

$startRow = 21
foreach ($sheet in @($ws1, $ws2)) {
  for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $pathCell = $sheet.Cells.Item($r, 1)
    $nameCell = $sheet.Cells.Item($r, 2)
    $pathCell.Value = $newRows[$i][0]
    $nameCell.NumberFormat = "@"
    $nameCell.Value = $newRows[$i][1]
  }
}

# Update selections on both sheets to cover the new data range
[void]$ws1.Activate()
[void]$ws1.Range("A2:B22").Select()

[void]$ws2.Activate()
[void]$ws2.Range("A2:B22").Select()

# Restore sheet1 as the active tab (it was tabSelected before the edit)
[void]$ws1.Activate()
